# Updates national_sdg_indicators.xlsx per the "Okno serwisowe 23.04.2024" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Number-format codes used by the workbook for the affected numeric columns
# (1 decimal / 2 decimals / 0 decimals, matching styles s="4", s="5", s="6").
$fmt1 = "[$-10809]0.0;\-0.0;0.0"
$fmt2 = "[$-10809]0.00;\-0.00;0.00"
$fmt0 = "[$-10809]0;\-0;0"

function Set-Cell {
    param(
        [string]$addr,
        $value,
        $fmt
    )
    $rng = $ws.Range($addr)
    if ($fmt -ne $null) {
        $rng.NumberFormat = $fmt
    }
    $rng.Value = $value
}

Set-Cell "Q6" 25.4 $null
Set-Cell "S6" 33.200000000000003 $null
Set-Cell "T6" 36 $fmt1
Set-Cell "T7" 51.4 $fmt1
Set-Cell "T8" 15.1 $fmt1
Set-Cell "R12" 0.46 $null
Set-Cell "S12" 0.82 $fmt2
Set-Cell "S19" 13.83 $fmt2
Set-Cell "T19" 14.74 $fmt2
Set-Cell "S26" 426.2 $fmt1
Set-Cell "S27" 406.2 $fmt1
Set-Cell "S28" 444.9 $fmt1
Set-Cell "S29" 430.6 $fmt1
Set-Cell "S30" 419.7 $fmt1
Set-Cell "S31" 253.9 $fmt1
Set-Cell "S32" 283.10000000000002 $fmt1
Set-Cell "S33" 226.6 $fmt1
Set-Cell "S34" 302 $fmt1
Set-Cell "S35" 222.2 $fmt1
Set-Cell "S36" 28.8 $fmt1
Set-Cell "S37" 27.2 $fmt1
Set-Cell "S38" 30.3 $fmt1
Set-Cell "S39" 29.2 $fmt1
Set-Cell "S40" 28.2 $fmt1
Set-Cell "S41" 26.6 $fmt1
Set-Cell "S42" 31.3 $fmt1
Set-Cell "S43" 22.2 $fmt1
Set-Cell "S44" 27.9 $fmt1
Set-Cell "S45" 24.6 $fmt1
Set-Cell "K65" 23.3 $fmt1
Set-Cell "M65" 35.299999999999997 $fmt1
Set-Cell "O65" 56.6 $fmt1
Set-Cell "R65" 68.400000000000006 $fmt1
Set-Cell "S65" 89 $fmt1
Set-Cell "T65" 83.3 $fmt1
Set-Cell "S118" 46.1 $fmt1
Set-Cell "S119" 33.1 $fmt1
Set-Cell "S120" 42.4 $fmt1
Set-Cell "S121" 32.299999999999997 $fmt1
Set-Cell "S122" 56.3 $fmt1
Set-Cell "S123" 35.4 $fmt1
Set-Cell "S124" 59.7 $fmt1
Set-Cell "Q125" 25.4 $null
Set-Cell "S125" 33.200000000000003 $null
Set-Cell "T125" 36 $fmt1
Set-Cell "T126" 51.4 $fmt1
Set-Cell "T127" 15.1 $fmt1
Set-Cell "R138" 15.61 $null
Set-Cell "S138" 16.88 $fmt2
Set-Cell "S141" 9.1 $fmt1
Set-Cell "S187" 44.2 $fmt1
Set-Cell "T188" 15.8 $fmt1
Set-Cell "T204" 170 $fmt0
Set-Cell "T205" 9.3000000000000007 $fmt1
Set-Cell "T207" 6.1 $fmt1
Set-Cell "S216" 14.8 $fmt1
Set-Cell "S222" 7.6 $fmt1
Set-Cell "S226" 64070 $fmt0
Set-Cell "N235" 100.7 $null
Set-Cell "S235" 94.2 $fmt1
Set-Cell "L236" 94.2 $null
Set-Cell "S236" 93.4 $fmt1
Set-Cell "R237" 15.61 $null
Set-Cell "S237" 16.88 $fmt2
Set-Cell "R238" 1189 $fmt0
Set-Cell "S238" 1318 $fmt0
Set-Cell "T239" 28779 $fmt0
Set-Cell "T247" 75 $fmt0
Set-Cell "T248" 7 $fmt0
Set-Cell "S249" 80 $fmt0
Set-Cell "T249" 97 $fmt0
Set-Cell "T289" 2603.83 $fmt2
Set-Cell "T290" 1744.48 $fmt2
Set-Cell "T294" 8 $fmt0
Set-Cell "T295" 0 $fmt0
Set-Cell "T296" 60 $fmt0

# Update the "Last update" footer text (row 300) shown under the data table.
$ws.Range("A300").Value = "Last update: 23-04-2024, 13:00"

# The 13.2.a renewable-energy row now cites Eurostat instead of Statistics Poland.
$ws.Range("U237").Value = "Eurostat"

# Refresh the saved selection so cell A1 (merged A1:B1, the title banner) is selected again.
$ws.Range("A1:B1").Select()

# Slightly widen the print scale used for the PDF/paper export.
$ws.PageSetup.Zoom = 62
